# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet between "2021-Q4" and "总计" with the
#    quarter's per-fund holdings.
# 2. Add a new leading row to the "总计" (totals) sheet summarising 2022-Q1.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the value to be stored as text (not auto-coerced to a number),
    # then drop back to the default "Normal" style so no stray number
    # format sticks around on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. New sheet "2022-Q1", positioned right after "2021-Q4"
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# Reuse the exact header formatting (bold/border/center) already used by the
# "2021-Q4" sheet header row so the new sheet's styling matches.
$afterSheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$afterSheet.Range("A2").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 - 创金合信港股通成长股票型发起式证券投资基金A
$q1.Cells.Item(2, 1).Value = 0
Set-TextValue $q1.Range("B2") "012315"
Set-TextValue $q1.Range("C2") "创金合信港股通成长股票型发起式证券投资基金A"
Set-TextValue $q1.Range("D2") "0.19"
Set-TextValue $q1.Range("E2") "83.49"
Set-TextValue $q1.Range("F2") "7.41"
Set-TextValue $q1.Range("G2") "0.0141"
$q1.Cells.Item(2, 8).Value = 4

# Row 3 - 创金合信港股通成长股票型发起式证券投资基金C
$q1.Cells.Item(3, 1).Value = 1
Set-TextValue $q1.Range("B3") "012316"
Set-TextValue $q1.Range("C3") "创金合信港股通成长股票型发起式证券投资基金C"
Set-TextValue $q1.Range("D3") "0.10"
Set-TextValue $q1.Range("E3") "83.49"
Set-TextValue $q1.Range("F3") "7.41"
Set-TextValue $q1.Range("G3") "0.0074"
$q1.Cells.Item(3, 8).Value = 4

# ---------------------------------------------------------------------------
# 2. "总计" sheet gains a new leading data row for 2022-Q1; existing rows
#    shift down one.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Copy the row-index cell format down onto the newly used A4 cell before
# rewriting the values (keeps the same bold/border/center formatting as
# A2 and A3).
$zj.Range("A3").Copy()
$zj.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

# Row 4 (was row 3): 2021-Q3
$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(4, 2).Value = "2021-Q3"
$zj.Cells.Item(4, 3).Value = 2
$zj.Cells.Item(4, 4).Value = 0.02

# Row 3 (was row 2): 2021-Q4
$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2021-Q4"
$zj.Cells.Item(3, 3).Value = 4
$zj.Cells.Item(3, 4).Value = 0.22

# Row 2 (new): 2022-Q1
$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q1"
$zj.Cells.Item(2, 3).Value = 2
$zj.Cells.Item(2, 4).Value = 0.02
